$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6276753333333334
$ws.Range("H2").Value = 1.883026
$ws.Range("I2").Value = 0.01507055680360577
$ws.Range("J2").Value = 0.01507055680360577
$ws.Range("M2").Value = 0.140567
$ws.Range("N2").Value = 0.421701
$ws.Range("O2").Value = 0.07810038533383065
$ws.Range("P2").Value = 0.07810038533383065
$ws.Range("Q2").Value = 0.08823043858066668
$ws.Range("R2").Value = 0.794073947226
$ws.Range("S2").Value = 0.001177016293556994
$ws.Range("T2").Value = 0.001177016293556994
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6276753333333334
$ws.Range("H3").Value = 1.883026
$ws.Range("I3").Value = 0.01507055680360577
$ws.Range("J3").Value = 0.01507055680360577
$ws.Range("O3").Value = 0.02984383293631935
$ws.Range("P3").Value = 0.02984383293631935
$ws.Range("Q3").Value = 0.03371474362955556
$ws.Range("R3").Value = 0.303432692666
$ws.Range("S3").Value = 0.0004497631795041215
$ws.Range("T3").Value = 0.0004497631795041215
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6276753333333334
$ws.Range("H4").Value = 1.883026
$ws.Range("I4").Value = 0.01507055680360577
$ws.Range("J4").Value = 0.01507055680360577
$ws.Range("O4").Value = 0.8920557817298499
$ws.Range("P4").Value = 0.8920557817298499
$ws.Range("Q4").Value = 1.007760365381333
$ws.Range("R4").Value = 9.069843288432001
$ws.Range("S4").Value = 0.01344377733054465
$ws.Range("T4").Value = 0.01344377733054465
$ws.Range("I5").Value = 0.6396010460118555
$ws.Range("J5").Value = 0.6396010460118555
$ws.Range("M5").Value = 0.140567
$ws.Range("N5").Value = 0.421701
$ws.Range("O5").Value = 0.07810038533383065
$ws.Range("P5").Value = 0.07810038533383065
$ws.Range("Q5").Value = 3.744538542383333
$ws.Range("R5").Value = 33.70084688145
$ws.Range("S5").Value = 0.04995308815344707
$ws.Range("T5").Value = 0.04995308815344707
$ws.Range("I6").Value = 0.6396010460118555
$ws.Range("J6").Value = 0.6396010460118555
$ws.Range("O6").Value = 0.02984383293631935
$ws.Range("P6").Value = 0.02984383293631935
$ws.Range("S6").Value = 0.01908814676307292
$ws.Range("T6").Value = 0.01908814676307292
$ws.Range("I7").Value = 0.6396010460118555
$ws.Range("J7").Value = 0.6396010460118555
$ws.Range("O7").Value = 0.8920557817298499
$ws.Range("P7").Value = 0.8920557817298499
$ws.Range("S7").Value = 0.5705598110953355
$ws.Range("T7").Value = 0.5705598110953355
$ws.Range("I8").Value = 0.3453283971845387
$ws.Range("J8").Value = 0.3453283971845388
$ws.Range("M8").Value = 0.140567
$ws.Range("N8").Value = 0.421701
$ws.Range("O8").Value = 0.07810038533383065
$ws.Range("P8").Value = 0.07810038533383065
$ws.Range("Q8").Value = 2.021721979818333
$ws.Range("R8").Value = 18.195497818365
$ws.Range("S8").Value = 0.02697028088682659
$ws.Range("T8").Value = 0.0269702808868266
$ws.Range("I9").Value = 0.3453283971845387
$ws.Range("J9").Value = 0.3453283971845388
$ws.Range("O9").Value = 0.02984383293631935
$ws.Range("P9").Value = 0.02984383293631935
$ws.Range("Q9").Value = 0.7725433459961112
$ws.Range("R9").Value = 6.952890113965001
$ws.Range("S9").Value = 0.01030592299374231
$ws.Range("T9").Value = 0.01030592299374231
$ws.Range("I10").Value = 0.3453283971845387
$ws.Range("J10").Value = 0.3453283971845388
$ws.Range("O10").Value = 0.8920557817298499
$ws.Range("P10").Value = 0.8920557817298499
$ws.Range("S10").Value = 0.3080521933039698
$ws.Range("T10").Value = 0.3080521933039698
